# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1:AF1 - match the look of the existing header row
# (bold, centered, thin border) by copying the format from AC1 first.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-45: every player row gets the team's season record.
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 98   # AD
    $ws.Cells.Item($r, 31).Value = 64   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
